$d = $word.ActiveDocument

# The "COMPETENCES TECHNIQUES" section lists 7 skill-category paragraphs.
# This edit reorders their text content (paragraph formatting/order of
# paragraph shells stays the same; only which text sits in which
# paragraph changes), going from:
#   1. Web : api, json
#   2. Langages : scala, java, python, matlab, c, c++
#   3. Bases de données : SQL, MongoDB, Neo4j, Redis
#   4. Autres :  c storage , bigquery, talend, mise en œuvre, cloud run, sub
#   5. Visualisation : tableau
#   6. ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn
#   7. MLOps : devops, nosql, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit
# to:
#   1. Langages : scala, java, python, matlab, c, c++
#   2. Visualisation : tableau
#   3. MLOps : devops, nosql, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit
#   4. Web : api, json
#   5. Autres :  c storage , bigquery, talend, mise en œuvre, cloud run, sub   (unchanged)
#   6. ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn (unchanged)
#   7. Bases de données : SQL, MongoDB, Neo4j, Redis

$newTexts = @(
  "Langages : scala, java, python, matlab, c, c++",
  "Visualisation : tableau",
  "MLOps : devops, nosql, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit",
  "Web : api, json",
  "Autres :  c storage , bigquery, talend, mise en œuvre, cloud run, sub",
  "ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn",
  "Bases de données : SQL, MongoDB, Neo4j, Redis"
)

# Locate the first of the 7 consecutive skills paragraphs by finding the
# "COMPETENCES TECHNIQUES" heading paragraph, the skills list starts right
# after it.
$startIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.TrimEnd([char]13, [char]7) -eq "COMPETENCES TECHNIQUES") {
        $startIndex = $i + 1
        break
    }
}

if ($startIndex -eq -1) {
    throw "Could not locate 'COMPETENCES TECHNIQUES' heading paragraph."
}

for ($j = 0; $j -lt $newTexts.Length; $j++) {
    $p = $d.Paragraphs.Item($startIndex + $j)
    $r = $d.Range($p.Range.Start, $p.Range.End - 1)
    $r.Text = $newTexts[$j]
}
